# Update "想去人数" (interested-people count) values in column F
# for two sheets: "展览" (sheet index 1) and "全部类型" (sheet index 4).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 94
$ws1.Range("F4").Value = 292
$ws1.Range("F6").Value = 119
$ws1.Range("F7").Value = 295
$ws1.Range("F9").Value = 2045
$ws1.Range("F11").Value = 4897
$ws1.Range("F12").Value = 95
$ws1.Range("F13").Value = 341

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 94
$ws4.Range("F6").Value = 292
$ws4.Range("F8").Value = 119
$ws4.Range("F9").Value = 295
$ws4.Range("F13").Value = 2045
$ws4.Range("F15").Value = 4897
$ws4.Range("F16").Value = 95
$ws4.Range("F17").Value = 341
